# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

# OFF sheet - update Home (row 2) target depth data
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 286
$wsOff.Range("C2").Value = 219
$wsOff.Range("D2").Value = 71
$wsOff.Range("E2").Value = 43

# DEF sheet - update Home (row 2) target depth data
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 321
$wsDef.Range("C2").Value = 224
$wsDef.Range("D2").Value = 68
$wsDef.Range("E2").Value = 31
$wsDef.Range("F2").Value = 6
